$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 133, shifting existing rows 133-135 down to 134-136.
$ws.Rows.Item(133).Insert()

# Fill the new row 133 with the full record, copying the static fields from the
# row that was pushed down to 134, and setting the new/changed values.
$ws.Cells.Item(133, 1).Value = 10
$ws.Cells.Item(133, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(133, 3).Value = "La Araucanía"
$ws.Cells.Item(133, 4).Value = 44448
$ws.Cells.Item(133, 4).NumberFormat = $ws.Cells.Item(134, 4).NumberFormat
$ws.Cells.Item(133, 5).Value = 9
$ws.Cells.Item(133, 6).Value = 100112039
$ws.Cells.Item(133, 7).Value = "Ciboulette"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 65
$ws.Cells.Item(133, 11).Value = 8000
$ws.Cells.Item(133, 12).Value = 8000
$ws.Cells.Item(133, 13).Value = 8000
$ws.Cells.Item(133, 14).Value = "`$/docena de atados"
$ws.Cells.Item(133, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(133, 16).Value = 2667
$ws.Cells.Item(133, 17).Value = 3
$ws.Cells.Item(133, 18).Value = "Hortaliza"
